# Auto-generated edit script applying numeric corrections to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 505.66666
$ws.Range("I32").Value = 342.85715
$ws.Range("J32").Value = 733.6
$ws.Range("K32").Value = 342.85715
$ws.Range("L32").Value = 733.6
$ws.Range("M32").Value = -16.85714999999999
$ws.Range("N32").Value = -1385.6
$ws.Range("H40").Value = 1198.2667
$ws.Range("I40").Value = 1053.4166
$ws.Range("K40").Value = 1053.4166
$ws.Range("M40").Value = -878.4166
$ws.Range("H64").Value = 4123.048
$ws.Range("I64").Value = 3959
$ws.Range("J64").Value = 4533.1665
$ws.Range("K64").Value = 3959
$ws.Range("L64").Value = 4533.1665
$ws.Range("M64").Value = -3711
$ws.Range("N64").Value = -5029.1665
$ws.Range("H67").Value = 4123.048
$ws.Range("I67").Value = 3959
$ws.Range("J67").Value = 4533.1665
$ws.Range("K67").Value = 3959
$ws.Range("L67").Value = 4533.1665
$ws.Range("M67").Value = -3101
$ws.Range("N67").Value = -6249.1665
$ws.Range("H100").Value = 2855.7144
$ws.Range("I100").Value = 2098
$ws.Range("K100").Value = 2098
$ws.Range("M100").Value = -1557
$ws.Range("H129").Value = 840.24
$ws.Range("I129").Value = 800
$ws.Range("J129").Value = 840.6464999999999
$ws.Range("K129").Value = 2400
$ws.Range("L129").Value = 2521.9395
$ws.Range("M129").Value = 2600
$ws.Range("N129").Value = -12521.9395
$ws.Range("H132").Value = 4792.5625
$ws.Range("I132").Value = 4978.3335
$ws.Range("J132").Value = 2006
$ws.Range("K132").Value = 14935.0005
$ws.Range("L132").Value = 6018
$ws.Range("M132").Value = -12405.0005
$ws.Range("N132").Value = -11078
$ws.Range("H138").Value = 1676.3871
$ws.Range("I138").Value = 558.1111
$ws.Range("J138").Value = 3224.7693
$ws.Range("K138").Value = 1674.3333
$ws.Range("L138").Value = 9674.3079
$ws.Range("M138").Value = 3465.6667
$ws.Range("N138").Value = -19954.3079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1594.0646
$ws.Range("I2").Value = 1029.8182
$ws.Range("J2").Value = 2973.3333
$ws.Range("K2").Value = 1029.8182
$ws.Range("L2").Value = 2973.3333
$ws.Range("M2").Value = -916.8181999999999
$ws.Range("N2").Value = -3199.3333
$ws.Range("H32").Value = 2685.5557
$ws.Range("I32").Value = 2082.5652
$ws.Range("K32").Value = 2082.5652
$ws.Range("M32").Value = -1795.5652
$ws.Range("H43").Value = 28781.75
$ws.Range("I43").Value = 30375
$ws.Range("J43").Value = 27188.5
$ws.Range("K43").Value = 30375
$ws.Range("L43").Value = 27188.5
$ws.Range("M43").Value = -30062
$ws.Range("N43").Value = -27814.5
$ws.Range("H107").Value = 50000
$ws.Range("J107").Value = 50000
$ws.Range("L107").Value = 50000
$ws.Range("N107").Value = -57680
$ws.Range("H110").Value = 1884.2941
$ws.Range("I110").Value = 2525.125
$ws.Range("J110").Value = 1314.6666
$ws.Range("K110").Value = 2525.125
$ws.Range("L110").Value = 1314.6666
$ws.Range("M110").Value = -480.125
$ws.Range("N110").Value = -5404.6666
$ws.Range("H116").Value = 1594.0646
$ws.Range("I116").Value = 1029.8182
$ws.Range("J116").Value = 2973.3333
$ws.Range("K116").Value = 1029.8182
$ws.Range("L116").Value = 2973.3333
$ws.Range("M116").Value = 1264.1818
$ws.Range("N116").Value = -7561.3333
$ws.Range("H122").Value = 1706.0667
$ws.Range("I122").Value = 1749.2
$ws.Range("J122").Value = 1619.8
$ws.Range("K122").Value = 5247.6
$ws.Range("L122").Value = 4859.4
$ws.Range("M122").Value = -2797.6
$ws.Range("N122").Value = -9759.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1594.0646
$ws.Range("I3").Value = 1029.8182
$ws.Range("J3").Value = 2973.3333
$ws.Range("K3").Value = 1029.8182
$ws.Range("L3").Value = 2973.3333
$ws.Range("M3").Value = -915.8181999999999
$ws.Range("N3").Value = -3201.3333
$ws.Range("H64").Value = 45455220
$ws.Range("I64").Value = 142858460
$ws.Range("J64").Value = 367.33334
$ws.Range("K64").Value = 142858460
$ws.Range("L64").Value = 367.33334
$ws.Range("M64").Value = -142858235
$ws.Range("N64").Value = -817.33334
$ws.Range("H67").Value = 45455220
$ws.Range("I67").Value = 142858460
$ws.Range("J67").Value = 367.33334
$ws.Range("K67").Value = 142858460
$ws.Range("L67").Value = 367.33334
$ws.Range("M67").Value = -142857680
$ws.Range("N67").Value = -1927.33334
$ws.Range("H81").Value = 20136.166
$ws.Range("J81").Value = 20136.166
$ws.Range("L81").Value = 20136.166
$ws.Range("N81").Value = -22258.166
$ws.Range("H84").Value = 20136.166
$ws.Range("J84").Value = 20136.166
$ws.Range("L84").Value = 60408.49800000001
$ws.Range("N84").Value = -71016.49800000001
$ws.Range("H99").Value = 2555.4285
$ws.Range("I99").Value = 2151.3333
$ws.Range("K99").Value = 2151.3333
$ws.Range("M99").Value = -653.3332999999998
$ws.Range("H137").Value = 50765
$ws.Range("J137").Value = 50765
$ws.Range("L137").Value = 50765
$ws.Range("N137").Value = -60965

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2428.8333
$ws.Range("I107").Value = 1712
$ws.Range("J107").Value = 6013
$ws.Range("K107").Value = 1712
$ws.Range("L107").Value = 6013
$ws.Range("M107").Value = 208
$ws.Range("N107").Value = -9853

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2000
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 6000
$ws.Range("N80").Value = -7872
$ws.Range("H83").Value = 2000
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27360
$ws.Range("H110").Value = 5300.25
$ws.Range("I110").Value = 1100.5
$ws.Range("J110").Value = 9500
$ws.Range("K110").Value = 3301.5
$ws.Range("L110").Value = 28500
$ws.Range("M110").Value = 788.5
$ws.Range("N110").Value = -36680

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1673.7916
$ws.Range("I102").Value = 1708.5883
$ws.Range("K102").Value = 1708.5883
$ws.Range("M102").Value = -86.58829999999989
$ws.Range("H132").Value = 61364.223
$ws.Range("I132").Value = 7820
$ws.Range("J132").Value = 104199.6
$ws.Range("K132").Value = 23460
$ws.Range("L132").Value = 312598.8
$ws.Range("M132").Value = -20930
$ws.Range("N132").Value = -317658.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2565
$ws.Range("I46").Value = 2463.3333
$ws.Range("J46").Value = 2666.6667
$ws.Range("K46").Value = 2463.3333
$ws.Range("L46").Value = 2666.6667
$ws.Range("M46").Value = -2275.3333
$ws.Range("N46").Value = -3042.6667
$ws.Range("H136").Value = 500000.34
$ws.Range("I136").Value = 500000.34
$ws.Range("K136").Value = 1500001.02
$ws.Range("M136").Value = -1497451.02
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 17250
$ws.Range("J54").Value = 17250
$ws.Range("L54").Value = 17250
$ws.Range("N54").Value = -18290
$ws.Range("H96").Value = 7093.2
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 7093.2
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 7093.2
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -9839.200000000001
$ws.Range("H107").Value = 3031433.5
$ws.Range("I107").Value = 367
$ws.Range("J107").Value = 5683616.5
$ws.Range("K107").Value = 1101
$ws.Range("L107").Value = 17050849.5
$ws.Range("M107").Value = 819
$ws.Range("N107").Value = -17054689.5
$ws.Range("H126").Value = 1614.4166
$ws.Range("I126").Value = 996.625
$ws.Range("J126").Value = 2850
$ws.Range("K126").Value = 2989.875
$ws.Range("L126").Value = 8550
$ws.Range("M126").Value = -519.875
$ws.Range("N126").Value = -13490
$ws.Range("H132").Value = 2800.9167
$ws.Range("I132").Value = 2401.5293
$ws.Range("K132").Value = 7204.5879
$ws.Range("M132").Value = -4674.5879

Write-Host "Applied all cell updates"
